# The edit adds a timestamp (date the source data file was refreshed) to
# cell C1 of the "About" sheet, formatted as a short date.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("About")

# Set the number format before the value so the new style slot picks up the
# built-in short-date format (numFmtId 14) instead of minting a custom one.
$ws.Range("C1").NumberFormat = "mm-dd-yy"
$ws.Range("C1").Value = Get-Date -Year 2021 -Month 4 -Day 21 -Hour 0 -Minute 0 -Second 0
